# Iteration 2: populate the "OMDB Queries" worksheet with the team's
# planned queries (No. / Covers Sections / Query in plain English),
# and update the selected cell to reflect where editing left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OMDB Queries")

# Data for rows 2-6: column B = sections covered, column C = query text.
# Rows 2-4 were filled in B-then-C order; rows 5-6 had their query text
# (column C) typed first and the covered-sections (column B) filled in
# afterwards - mirrored here so shared-string insertion order matches.
function Set-QueryText($row, $text) {
    $c = $ws.Cells.Item($row, 3)
    $c.Value = $text
    $c.Font.Name = "Arial"
    $c.Font.Size = 12
    $c.Font.Color = 1114146
    $c.HorizontalAlignment = -4131
    $c.VerticalAlignment = -4108
    $c.ReadingOrder = 1
}

$ws.Cells.Item(2, 2).Value = "A + C"
Set-QueryText 2 "Display movies with select actor(s)"

$ws.Cells.Item(3, 2).Value = "A"
Set-QueryText 3 "Display movies made after 1990 "

$ws.Cells.Item(4, 2).Value = "A + D"
Set-QueryText 4 "Find movies with select song(s)"

Set-QueryText 5 "Display movies of a select genre"
Set-QueryText 6 "Display movies with select actor, and select genre"

$ws.Cells.Item(5, 2).Value = "A + B"
$ws.Cells.Item(6, 2).Value = "A + B + C"

# Leave the selection where the author's cursor ended up after data entry.
$ws.Range("B7").Select()
